$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update shared text: "Thomas Hex" -> "Matthies Hex" (row 9, column B)
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "Matthies Hex"

# ---------------------------------------------------------------------------
# 2. Rebuild the numeric simulation block (columns C:W, rows 4-31).
#    Row 3 is left untouched. Two brand new rows of simulated data were
#    inserted right after row 3 (now rows 4 and 5), pushing the previous
#    rows 4-29 down to rows 6-31. Two more rows of source data (28, 29)
#    were appended at the bottom (now rows 30 and 31).
# ---------------------------------------------------------------------------
$blob = "0.9767157330826673,1.005840657983691,1.01462516818683,0.9922005644938888,1.014033140182033,1.014033140182033,1.014033140182033,0.9627576361178829,1.015672225810816,0.9964229459858126,0.9627576361178829,1.014033140182033,1.015672225810816,0.9892149309643496,1.003936395152353,0.9974876673702441,0.9902101421408628,0.9974876673702441,0.9961658916511553,0.999739341357331,0.9972835089804528|0.9650923920813844,1.024507468359565,1.015426116463245,0.9879868070216542,1.073221630243332,1.073221630243332,1.073221630243332,0.9532306849605919,0.9890126615316255,0.9754613632854079,0.9532306849605919,1.073221630243332,0.9890126615316255,0.9711216732461088,0.9884997342766398,1.005154992245183,0.976743384504624,1.005154992245183,1.000862945939301,1.015334682800107,0.9979923904933508|0.9861850144372715,1.002433409195066,1.009122992404097,0.995324880391146,1.00509893231041,1.00509893231041,1.00509893231041,0.977456201485069,1.011519720322913,0.9990508776790926,0.977456201485069,1.00509893231041,1.011519720322913,0.9944879609039912,1.00342230035703,0.9980249513727975,0.9947669340663761,0.9980249513727975,0.9973499336273846,0.9988997333639895,0.9982740035281331|0.975546426304032,1.006090920504322,1.015392836311245,0.991756186340055,1.014745029272334,1.014745029272334,1.014745029272334,0.9609771883285302,1.016523729538905,0.9962307585230558,0.9609771883285302,1.014745029272334,1.016523729538905,0.9887504589337177,1.00413995793948,0.9974153157132566,0.9897523680691634,0.9974153157132566,0.9960005333699562,0.9997494325504318,0.9971578843903099|0.9150356988405395,1.015522304936957,1.055912357731568,0.9711174920937969,1.033533966906653,1.033533966906653,1.033533966906653,0.8619257604103421,1.06956396383293,0.9933416639790058,0.8619257604103421,1.033533966906653,1.06956396383293,0.9657448621216358,1.020340727963363,0.9883412303833085,0.9675357387790228,0.9883412303833085,0.9840352958109305,0.9939350300300751,0.9894941510914741|0.9977945964751253,1.001279195195088,1.001086124496262,0.999244895817993,1.003738843579841,1.003738843579841,1.003738843579841,0.9968937483569253,0.9998941666680996,0.9987754043556544,0.9968937483569253,1.003738843579841,0.9998941666680996,0.9983939575125125,0.9995695312430464,1.000175586201622,0.9986776036143393,1.000175586201622,0.9999429136057146,1.00070209960054,0.9998383718681235|0.9998692573074208,1.000135822963761,1.000041234313033,0.9999484644368877,1.000435916883597,1.000435916883597,1.000435916883597,0.9998622518457062,0.9998593098681134,0.9998476348214227,0.9998622518457062,1.000435916883597,0.9998593098681134,0.9998607808569098,0.9999038871525006,1.000052492865805,0.9998900087169025,1.000052492865805,1.000026485758576,1.00010837198358,0.9999999865549929|0.9963299276036417,1.002110686467543,1.001815313405847,0.9987423385640173,1.006166108757411,1.006166108757411,1.006166108757411,0.9948234880996217,0.9998626105243392,0.9979823566857609,0.9948234880996217,1.006166108757411,0.9998626105243392,0.9973430493119805,0.9993024745441783,1.000284069127124,0.9978094790626594,1.000284069127124,0.9998986364863472,1.00115213094056,0.9997291037635228|0.9136018831670429,1.015519000956498,1.056978506134661,0.9705898976380296,1.033348530123681,1.033348530123681,1.033348530123681,0.8595405382161821,1.071293254226499,0.9934969588170566,0.8595405382161821,1.033348530123681,1.071293254226499,0.9654168962213407,1.020941575932264,0.9880607741887873,0.9671412300269037,0.9880607741887873,0.9836930550510979,0.9936241500656144,0.9892960711599563|0.9167739046534271,1.029291852469156,1.04866936325695,0.9724138699283404,1.076744816608904,1.076744816608904,1.076744816608904,0.8707353459223337,1.037852239602815,0.9775306339027066,0.8707353459223337,1.076744816608904,1.037852239602815,0.9542937927625744,1.005133054765578,0.9951108007113509,0.9603338184844965,0.9951108007113509,0.9894365680155983,1.006898217734259,0.9912515032930791|0.9713837227894748,1.005758654336842,1.018604709515791,0.9902884533052621,1.012977091894737,1.012977091894737,1.012977091894737,0.9537437993894713,1.022282591252629,0.9971440846210531,0.9537437993894713,1.012977091894737,1.022282591252629,0.9880131953210504,1.006285522278946,0.9963344941789458,0.9887716146491209,0.9963344941789458,0.9948229839605249,0.9984538055473674,0.9965228883881576|1.01522048619772,0.994530066122603,0.9911854778758697,1.004910080810938,0.9859740104856882,0.9859740104856882,0.9859740104856882,1.023868500365388,0.9932486765746102,1.004080440597878,1.023868500365388,0.9859740104856882,0.9932486765746102,1.008558588469999,0.9990793786927742,1.001030395808562,1.007342419250312,1.001030395808562,1.002000317059156,0.9987950557444625,1.001627217378837|0.8503200600000009,1.026284799999999,1.0989676,0.9490348699999991,1.055856,1.055856,1.055856,0.7563715000000001,1.124812899999999,0.9894340199999996,0.7563715000000001,1.055856,1.124812899999999,0.9405921999999998,1.036923884999999,0.9790134666666667,0.9434064233333329,0.9790134666666664,0.9715188174999996,0.9883862539999997,0.9813852187499998|1.5361801,0.81025805,0.687925,1.1741023,0.5123396,0.5123396,0.5123396,1.8399756,0.7563714999999999,1.1415979,1.8399756,0.5123396,0.7563714999999999,1.29817355,0.9652368999999998,1.0362289,1.256816466666667,1.0362289,1.07069725,0.9590257199999999,1.05734375625|0.85032006,1.0262848,1.0989676,0.94903487,1.055856,1.055856,1.055856,0.7563714999999999,1.1248129,0.98943402,0.7563714999999999,1.055856,1.1248129,0.9405922,1.036923885,0.9790134666666667,0.9434064233333332,0.9790134666666667,0.9715188175,0.9883862539999999,0.98138521875|1.3647435,0.8658522200000001,0.7893455499999998,1.1201163,0.64704416,0.64704416,0.64704416,1.5651105,0.84625205,1.1044132,1.5651105,0.64704416,0.84625205,1.205681275,0.9831841750000001,1.019468903333334,1.177159616666667,1.019468903333333,1.0446307525,0.965113434,1.037859685|0.9644909860273967,1.02283253739726,1.016588058630136,0.9877157665753427,1.067833660821918,1.067833660821918,1.067833660821918,0.9514463479452052,0.9933552183561645,0.9774623400000002,0.9514463479452052,1.067833660821918,0.9933552183561645,0.9724007831506849,0.9905354924657536,1.004211742374429,0.9775057776255709,1.004211742374429,1.000087748424658,1.01363693090411,0.9977156144691779|1.1752191,0.9249589584210527,0.903416825263158,1.057086286315789,0.7976590484210528,0.7976590484210528,0.7976590484210528,1.267154015789474,0.9488797357894737,1.062052425263158,1.267154015789474,0.7976590484210528,0.9488797357894737,1.108016875789474,1.002983011052632,1.004564266666667,1.091040012631579,1.004564266666667,1.017694771578947,0.9736876269473687,1.017053299407895|0.9690087131578947,1.02012635368421,1.014431427368421,0.9891503026315788,1.060209326842105,1.060209326842105,1.060209326842105,0.9580039736842104,0.9936964242105264,0.9799370626315789,0.9580039736842104,1.060209326842105,0.9936964242105264,0.9758501989473685,0.9914233634210525,1.003969908245614,0.9802835668421053,1.003969908245614,1.000265006842105,1.012253870842105,0.9980704480263157|1.353230360303588,0.8710206287832202,0.7956458042010981,1.11617974601238,0.6616320149711334,0.6616320149711334,0.6616320149711334,1.548079765550479,0.8489902940146983,1.099827434327867,1.548079765550479,0.6616320149711334,0.8489902940146983,1.198535029782589,0.9825850200135392,1.01956735817877,1.171083268525852,1.01956735817877,1.043720455137173,0.967302767103965,1.036825756020558|1.09107527314148,0.9623121884336756,0.9489378720293017,1.03074408250835,0.8961192061560384,0.8961192061560384,0.8961192061560384,1.137269315054173,0.9711561553077607,1.031982267584797,1.137269315054173,0.8961192061560384,0.9711561553077607,1.054212735180967,1.000950118908055,1.001514892172658,1.046389850956761,1.001514892172657,1.00882218975658,0.986281593036472,1.008699545026947|0.9969380340360081,1.002520573470641,1.001036825181439,0.9995049579417656,1.006068376289979,1.006068376289979,1.006068376289979,0.9948998809505066,0.9985390216161666,0.9980826254893002,0.9948998809505066,1.006068376289979,0.9985390216161666,0.9967194512833366,0.9990219897789661,0.9998357596188839,0.997647953502813,0.9998357596188839,0.9997530591996044,1.001016122617679,0.9996987868719758|0.9987088764496752,0.9995739654508983,1.001307491005544,0.998933926887475,1.000107663060356,1.000107663060356,1.000107663060356,0.9988911172967686,1.002150129736611,0.9999278953506071,0.9988911172967686,1.000107663060356,1.002150129736611,1.00052062351669,1.000542028312043,1.000382970031245,0.9999917246402849,1.000382970031245,1.000020709245303,1.000038100008313,0.9999501331547418|0.7386952690726108,1.118876613314101,1.14158724860542,0.9132417776256819,1.328982962628722,1.328982962628722,1.328982962628722,0.6087752629102897,1.060157332243438,0.8971795524067498,0.6087752629102897,1.328982962628722,1.060157332243438,0.8344662975768637,0.9866995549345599,0.9993051859274832,0.8607247909264698,0.9993051859274832,0.9777893338520328,1.048028059607371,0.9759370023508765|0.9219961428069141,1.022439094449975,1.047577256399734,0.9746200056586706,1.054272234625513,1.054272234625513,1.054272234625513,0.8751624151924474,1.046662981219421,0.9854981213623133,0.8751624151924474,1.054272234625513,1.046662981219421,0.9609126982059342,1.010641493439046,0.9920325436791272,0.9654818006901796,0.9920325436791271,0.987679409174013,1.000997974264313,0.9910285314643736|0.9926012164310961,1.003614436883144,1.003868471422279,0.9976702341828176,1.009769760072786,1.009769760072786,1.009769760072786,0.988783609356974,1.001232622467071,0.996945985930787,0.988783609356974,1.009769760072786,1.001232622467071,0.9950081159120223,0.9994514283249442,0.9999286639656101,0.9958954886689542,0.9999286639656101,0.999364056519912,1.001445197230487,0.9993107920933693|1.009885432499335,0.9994885432086125,0.9929710979543592,1.003297378827203,1.000489820470242,1.000489820470242,1.000489820470242,1.016888167018523,0.9890519015804752,0.9991541466020808,1.016888167018523,1.000489820470242,0.9890519015804752,1.002970034299499,0.9961746402038389,1.002143296356413,1.003079149142067,1.002143296356413,1.002431816974111,1.002043417673337,1.001403311020104|1.037923710137449,0.9908856899009137,0.976251528392305,1.011875123358037,0.9807563861481579,0.9807563861481579,0.9807563861481579,1.062640662584481,0.9731460358580927,1.004346705225237,1.062640662584481,0.9807563861481579,0.9731460358580927,1.017893349221287,0.9925105796080649,1.005514361530244,1.015887273933537,1.005514361530244,1.007104551987192,1.001834918819385,1.004728230200584"
$rowStrings = $blob -split '\|'

$numRows = $rowStrings.Count
$numCols = 21
$data = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $cellStrings = $rowStrings[$i] -split ','
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i, $j] = [double]$cellStrings[$j]
    }
}
$ws.Range("C4:W31").Value = $data

# ---------------------------------------------------------------------------
# 3. Add the row labels (columns A and B) for the two brand-new rows 30/31.
#    Column A keeps the bold/centered/bordered style used by the rest of
#    column A, so copy formatting down from the row above first.
# ---------------------------------------------------------------------------
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Michael-SNHex"
